# Sheet "Sibirev I. V." grade sheet: fill in a couple of missing homework
# scores and extend the totals formula in K to cover the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (student #6) ------------------------------------------------
# H9 was blank; the student got a 5 there too, so the total (K9) grows
# from 25 to 30 and now needs to explicitly span C9:H9.
$ws.Range("H9").Value = 5
$ws.Range("K9").Formula = "=SUM(C9:H9)"

# --- Row 20 (student #17) -----------------------------------------------
# G20 and H20 were blank; both scored a 5. Two more homework columns (I20,
# J20) are also recorded for this student - I20 left blank, J20 = 5 -
# matching the formatting already used for the same columns in row 4.
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 5

$ws.Range("I4").Copy()
$ws.Range("I20").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("J4").Copy()
$ws.Range("J20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J20").Value = 5

$excel.CutCopyMode = 0

# The total (K20) now spans the wider C20:J20 range, going from 20 to 35.
$ws.Range("K20").Formula = "=SUM(C20:J20)"

# Leave the selection on the last cell that was edited, H9.
$ws.Range("H9").Select()
